$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Model R²
$ws.Range("B7").Value = 0.347
$ws.Range("C7").Value = 0.674
$ws.Range("D7").Value = 0.543
$ws.Range("E7").Value = 0.652
$ws.Range("F7").Value = 0.505
$ws.Range("G7").Value = 0.552

# Row 8 - Model Adj R²
$ws.Range("B8").Value = -13.236
$ws.Range("C8").Value = 0.449
$ws.Range("D8").Value = -1.287
$ws.Range("E8").Value = 0.601
$ws.Range("F8").Value = 0.485
$ws.Range("G8").Value = 0.541

# Row 9 - Model RMSE
$ws.Range("B9").Value = 2.406
$ws.Range("C9").Value = 2.02
$ws.Range("D9").Value = 2.126
$ws.Range("E9").Value = 2.234
$ws.Range("F9").Value = 2.328
$ws.Range("G9").Value = 2.287

# Row 10 - Model HH
$ws.Range("C10").Value = 3
$ws.Range("E10").Value = 23
$ws.Range("F10").Value = 55
$ws.Range("G10").Value = 84

# Row 11 - Delta R²
$ws.Range("B11").Value = 1.016
$ws.Range("C11").Value = 0.36
$ws.Range("D11").Value = 0.6
$ws.Range("E11").Value = 0.654
$ws.Range("F11").Value = 0.537
$ws.Range("G11").Value = 0.5629999999999999

# Row 12 - Delta Adj R²
$ws.Range("B12").Value = 22.144
$ws.Range("C12").Value = 0.609
$ws.Range("D12").Value = 2.998
$ws.Range("E12").Value = 0.75
$ws.Range("F12").Value = 0.5580000000000001
$ws.Range("G12").Value = 0.578

# Row 13 - Delta RMSE
$ws.Range("B13").Value = -1.44
$ws.Range("C13").Value = -0.911
$ws.Range("D13").Value = -1.106
$ws.Range("E13").Value = -1.557
$ws.Range("F13").Value = -1.033
$ws.Range("G13").Value = -1.149

# Row 14 - Delta HH
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 10
